# "Fixed for summer dataset"
#
# Appointment Type Summation sheet: two new appointment types need to be
# added near the top of the list (right under the header row), pushing the
# existing list down. The row immediately below the header is left blank,
# and the two new appointment types land on the next two rows, with the
# rest of the original list shifted down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Appointment Type Summation")

# Shift the existing data rows (2-9) down by three rows (-> 5-12),
# working from the bottom up so we never overwrite a row before reading it.
$ws.Range("A12").Value = $ws.Range("A9").Text
$ws.Range("A11").Value = $ws.Range("A8").Text
$ws.Range("A10").Value = $ws.Range("A7").Text
$ws.Range("A9").Value  = $ws.Range("A6").Text
$ws.Range("A8").Value  = $ws.Range("A5").Text
$ws.Range("A7").Value  = $ws.Range("A4").Text
$ws.Range("A6").Value  = $ws.Range("A3").Text
$ws.Range("A5").Value  = $ws.Range("A2").Text

# Fill in the two new appointment types (row 2 stays blank).
$ws.Range("A4").Value = "Health & Science Coach: 30-Minute Meeting "
$ws.Range("A3").Value = "Business Coach: 30-Minute Meeting"
$ws.Range("A2").Value = ""
